# SA-HW10.xlsx — add "Gaussian-Quadrature" + 3 new spiral sampling-scheme rows
# (averaged-intensity results for spiral quadrature runs), reshuffling the
# existing "Gaussian-Quadrature" row to sit right after the "Ring Perpendicular"
# rows and pushing the rotation/hex-grid rows down to make room.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 17-19 need the same bold/bordered look as the rest of column A,
# so clone the formatting from the rows directly above before filling them in.
$ws.Range("A14:A16").Copy($ws.Range("A17:A19"))
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(19, 1).Value = 17

# Row 10: Gaussian-Quadrature
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(10, 3).Value = 0.9783622755493311
$ws.Cells.Item(10, 4).Value = 0.7742431438492958
$ws.Cells.Item(10, 5).Value = 0.8977956654838386
$ws.Cells.Item(10, 6).Value = 0.8651770259480333
$ws.Cells.Item(10, 7).Value = 0.9783622755493311
$ws.Cells.Item(10, 8).Value = 0.7742431438492958
$ws.Cells.Item(10, 9).Value = 1.17407315656008
$ws.Cells.Item(10, 10).Value = 0.8497568785102666
$ws.Cells.Item(10, 11).Value = 1.078894377174901
$ws.Cells.Item(10, 12).Value = 0.8964656254791129
$ws.Cells.Item(10, 13).Value = 0.9783622755493311
$ws.Cells.Item(10, 14).Value = 0.8360194046665672
$ws.Cells.Item(10, 15).Value = 0.8788945277076247
$ws.Cells.Item(10, 16).Value = 0.9393460185693574

# Row 11: Spiral-90deg-10rot-5space
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11, 3).Value = 1.248855457822235
$ws.Cells.Item(11, 4).Value = 0.07329703369301814
$ws.Cells.Item(11, 5).Value = 1.359483459703277
$ws.Cells.Item(11, 6).Value = 0.8738971374394071
$ws.Cells.Item(11, 7).Value = 1.248855457822235
$ws.Cells.Item(11, 8).Value = 0.07329703369301814
$ws.Cells.Item(11, 9).Value = 1.200564394583324
$ws.Cells.Item(11, 10).Value = 1.168875103966613
$ws.Cells.Item(11, 11).Value = 0.9300433111305988
$ws.Cells.Item(11, 12).Value = 0.4817158307967824
$ws.Cells.Item(11, 13).Value = 1.248855457822235
$ws.Cells.Item(11, 14).Value = 0.7163902466981478
$ws.Cells.Item(11, 15).Value = 0.8888832721644844
$ws.Cells.Item(11, 16).Value = 0.917091466141907

# Row 12: Spiral-90deg-15rot-5space
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12, 3).Value = 1.247774295371553
$ws.Cells.Item(12, 4).Value = 0.07342730292086726
$ws.Cells.Item(12, 5).Value = 1.363652708218439
$ws.Cells.Item(12, 6).Value = 0.8756396546041382
$ws.Cells.Item(12, 7).Value = 1.247774295371553
$ws.Cells.Item(12, 8).Value = 0.07342730292086726
$ws.Cells.Item(12, 9).Value = 1.198145852783202
$ws.Cells.Item(12, 10).Value = 1.171321402803311
$ws.Cells.Item(12, 11).Value = 0.9288338991225245
$ws.Cells.Item(12, 12).Value = 0.4827602126371013
$ws.Cells.Item(12, 13).Value = 1.247774295371553
$ws.Cells.Item(12, 14).Value = 0.7185400055696531
$ws.Cells.Item(12, 15).Value = 0.8901234902787494
$ws.Cells.Item(12, 16).Value = 0.917694416057642

# Row 13: Spiral-90deg-10rot-3space
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13, 3).Value = 1.247381320248102
$ws.Cells.Item(13, 4).Value = 0.07326047195728534
$ws.Cells.Item(13, 5).Value = 1.362378545401238
$ws.Cells.Item(13, 6).Value = 0.8747960804149979
$ws.Cells.Item(13, 7).Value = 1.247381320248102
$ws.Cells.Item(13, 8).Value = 0.07326047195728534
$ws.Cells.Item(13, 9).Value = 1.200521685994359
$ws.Cells.Item(13, 10).Value = 1.169207473477452
$ws.Cells.Item(13, 11).Value = 0.9300442729081312
$ws.Cells.Item(13, 12).Value = 0.4818118298024217
$ws.Cells.Item(13, 13).Value = 1.247381320248102
$ws.Cells.Item(13, 14).Value = 0.7178195086792619
$ws.Cells.Item(13, 15).Value = 0.8894541045054059
$ws.Cells.Item(13, 16).Value = 0.9174252100254985

# Row 14: NoRotation-tilt60deg
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(14, 3).Value = 0.3305079999999967
$ws.Cells.Item(14, 4).Value = 0.03816400000000007
$ws.Cells.Item(14, 5).Value = 1.162688
$ws.Cells.Item(14, 6).Value = 0.9542199999999997
$ws.Cells.Item(14, 7).Value = 0.3305079999999967
$ws.Cells.Item(14, 8).Value = 0.03816400000000007
$ws.Cells.Item(14, 9).Value = 1.437944000000002
$ws.Cells.Item(14, 10).Value = 0.8471760000000002
$ws.Cells.Item(14, 11).Value = 1.325852000000001
$ws.Cells.Item(14, 12).Value = 0.3745759999999998
$ws.Cells.Item(14, 13).Value = 0.3305079999999967
$ws.Cells.Item(14, 14).Value = 0.6004259999999999
$ws.Cells.Item(14, 15).Value = 0.621394999999999
$ws.Cells.Item(14, 16).Value = 0.8088909999999998

# Row 15: Rotation-NoTilt
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Cells.Item(15, 3).Value = 0.4095874999999999
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0.5927875000000016
$ws.Cells.Item(15, 6).Value = 0.8227750000000015
$ws.Cells.Item(15, 7).Value = 0.4095874999999999
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 1.577387499999999
$ws.Cells.Item(15, 10).Value = 0.5104125000000003
$ws.Cells.Item(15, 11).Value = 1.674025000000002
$ws.Cells.Item(15, 12).Value = 0.19
$ws.Cells.Item(15, 13).Value = 0.4095874999999999
$ws.Cells.Item(15, 14).Value = 0.2963937500000008
$ws.Cells.Item(15, 15).Value = 0.4562875000000007
$ws.Cells.Item(15, 16).Value = 0.7221218750000005

# Row 16: Rotation-60detTilt
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"
$ws.Cells.Item(16, 3).Value = 0.6949599760384022
$ws.Cells.Item(16, 4).Value = 0.396083935232
$ws.Cells.Item(16, 5).Value = 0.7861691326464018
$ws.Cells.Item(16, 6).Value = 0.9028008359935983
$ws.Cells.Item(16, 7).Value = 0.6949599760384022
$ws.Cells.Item(16, 8).Value = 0.396083935232
$ws.Cells.Item(16, 9).Value = 1.318388635443199
$ws.Cells.Item(16, 10).Value = 0.7275068732416011
$ws.Cells.Item(16, 11).Value = 1.375534365491198
$ws.Cells.Item(16, 12).Value = 0.5318745195520009
$ws.Cells.Item(16, 13).Value = 0.6949150515200011
$ws.Cells.Item(16, 14).Value = 0.5911265339392009
$ws.Cells.Item(16, 15).Value = 0.6950034699776005
$ws.Cells.Item(16, 16).Value = 0.8416647842048002

# Row 17: HexGrid-90degTilt5degRes
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17, 3).Value = 0.995526852596242
$ws.Cells.Item(17, 4).Value = 0.9957446280651487
$ws.Cells.Item(17, 5).Value = 0.9868460505376425
$ws.Cells.Item(17, 6).Value = 0.9919090989159894
$ws.Cells.Item(17, 7).Value = 0.995526852596242
$ws.Cells.Item(17, 8).Value = 0.9957446280651487
$ws.Cells.Item(17, 9).Value = 0.9942716407676722
$ws.Cells.Item(17, 10).Value = 0.9855713702367078
$ws.Cells.Item(17, 11).Value = 0.9925462142448677
$ws.Cells.Item(17, 12).Value = 0.9940568618051426
$ws.Cells.Item(17, 13).Value = 0.9954987230375346
$ws.Cells.Item(17, 14).Value = 0.9912953393013957
$ws.Cells.Item(17, 15).Value = 0.9925066575287557
$ws.Cells.Item(17, 16).Value = 0.9920590896461766

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18, 3).Value = 0.9013334036250382
$ws.Cells.Item(18, 4).Value = 1.033470172030369
$ws.Cells.Item(18, 5).Value = 0.9502916294589036
$ws.Cells.Item(18, 6).Value = 0.9386454422073174
$ws.Cells.Item(18, 7).Value = 0.9013334036250382
$ws.Cells.Item(18, 8).Value = 1.033470172030369
$ws.Cells.Item(18, 9).Value = 0.9380574701864866
$ws.Cells.Item(18, 10).Value = 1.092376736792457
$ws.Cells.Item(18, 11).Value = 0.9997066688280566
$ws.Cells.Item(18, 12).Value = 0.9617712624077653
$ws.Cells.Item(18, 13).Value = 0.9013334036250382
$ws.Cells.Item(18, 14).Value = 0.9918809007446361
$ws.Cells.Item(18, 15).Value = 0.9559351618304068
$ws.Cells.Item(18, 16).Value = 0.9769565981920492

# Row 19: HexGrid-60degTilt5degRes
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19, 3).Value = 0.9304766936225563
$ws.Cells.Item(19, 4).Value = 1.083185521174754
$ws.Cells.Item(19, 5).Value = 0.9555256776006366
$ws.Cells.Item(19, 6).Value = 1.026233510769793
$ws.Cells.Item(19, 7).Value = 0.9304766936225563
$ws.Cells.Item(19, 8).Value = 1.083185521174754
$ws.Cells.Item(19, 9).Value = 0.9291280846667651
$ws.Cells.Item(19, 10).Value = 1.032022717645343
$ws.Cells.Item(19, 11).Value = 0.9620864756394276
$ws.Cells.Item(19, 12).Value = 1.066398590217768
$ws.Cells.Item(19, 13).Value = 0.9304201384045221
$ws.Cells.Item(19, 14).Value = 1.019355599387695
$ws.Cells.Item(19, 15).Value = 0.9988553507919351
$ws.Cells.Item(19, 16).Value = 0.9981321589171304
